# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
#
# 1) "cumcontrol" -> "cumcontrol1" (same 6 lines, B/C text refreshed - diversified copy)
# 2) new sheet "cumcontrol2" inserted right before "dickpic" (a second diversified
#    variant of the same delay/sync/edge control lines, reusing the A-column keys)
# 3) "dickpic" and "boosters" sheets are left untouched content-wise; they simply
#    shift right in the tab order to make room for the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: rename cumcontrol -> cumcontrol1 and refresh its Text/Note columns
# ---------------------------------------------------------------------------
$cc1 = $wb.Worksheets.Item("cumcontrol")
$cc1.Name = "cumcontrol1"

$cc1.Cells.Item(2,2).Value = "trust me you wanna save it for this next one, it's the best one"
$cc1.Cells.Item(2,3).Value = "DELAY variant."

$cc1.Cells.Item(3,2).Value = "edge for me... what I'm about to send is gonna hit different"
$cc1.Cells.Item(3,3).Value = "DELAY. Send PPV."

$cc1.Cells.Item(4,2).Value = "I'm close too man, let's go at the same time... check this"
$cc1.Cells.Item(4,3).Value = "SYNC variant."

$cc1.Cells.Item(5,2).Value = "aight let's bust together... open this"
$cc1.Cells.Item(5,3).Value = "SYNC. Send PPV."

$cc1.Cells.Item(6,2).Value = "hold it... you're gonna last until I say otherwise"
$cc1.Cells.Item(6,3).Value = "EDGE variant."

$cc1.Cells.Item(7,2).Value = "bro don't you dare finish yet... I'm not even close to done"
$cc1.Cells.Item(7,3).Value = "CONTROL."

# ---------------------------------------------------------------------------
# Step 2: build cumcontrol2 as a copy of dickpic's sheet (so it inherits the
# exact same column widths / header+row styling), inserted right before
# dickpic, then overwrite its Name/Text/Note columns with the new variant.
# ---------------------------------------------------------------------------
$dp = $wb.Worksheets.Item("dickpic")
$dp.Copy($dp)
$cc2 = $wb.Worksheets.Item("dickpic (2)")
$cc2.Name = "cumcontrol2"

$cc2.Cells.Item(2,1).Value = "delay2"
$cc2.Cells.Item(2,2).Value = "edge just a bit more... the finale is worth it"
$cc2.Cells.Item(2,3).Value = "DELAY variant."

$cc2.Cells.Item(3,1).Value = "delay1"
$cc2.Cells.Item(3,2).Value = "don't finish yet bro... this last one? insane"
$cc2.Cells.Item(3,3).Value = "DELAY. Send PPV."

$cc2.Cells.Item(4,1).Value = "sync2"
$cc2.Cells.Item(4,2).Value = "I'm about to blow too, watch this and let's go"
$cc2.Cells.Item(4,3).Value = "SYNC variant."

$cc2.Cells.Item(5,1).Value = "sync1"
$cc2.Cells.Item(5,2).Value = "alright go time, let's finish this together... open it"
$cc2.Cells.Item(5,3).Value = "SYNC. Send PPV."

$cc2.Cells.Item(6,1).Value = "edge2"
$cc2.Cells.Item(6,2).Value = "hold it, I got one more thing to show you first"
$cc2.Cells.Item(6,3).Value = "EDGE variant."

$cc2.Cells.Item(7,1).Value = "edge1"
$cc2.Cells.Item(7,2).Value = "yo slow down... we're not done here"
$cc2.Cells.Item(7,3).Value = "CONTROL."

Write-Host "Final sheet order:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
